$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (HP / glucose)
$ws.Range("E10").Value = 0.5986548995531681
$ws.Range("F10").Value = 0.5987440318041106
$ws.Range("G10").Value = 0.5988283263768763

# Row 11 (HP / cornstover)
$ws.Range("D11").Value = 0.6028488348904371
$ws.Range("E11").Value = 0.5854047692975199
$ws.Range("F11").Value = 0.5857518754390745
$ws.Range("G11").Value = 0.586100821567949

# Row 12 (HP / sugarcane)
$ws.Range("D12").Value = 0.6027367034234539
$ws.Range("E12").Value = 0.591779165550224
$ws.Range("F12").Value = 0.5923570483142888
$ws.Range("G12").Value = 0.5924496053209782

# Row 13 (HP / corn)
$ws.Range("D13").Value = 0.6032563360714954
$ws.Range("E13").Value = 0.5927185846790763
$ws.Range("F13").Value = 0.5928064600981822
$ws.Range("G13").Value = 0.5928909045560399

# Row 14 (HP_neutral / glucose)
$ws.Range("E14").Value = 0.6215865611896321
$ws.Range("F14").Value = 0.6216761342309467
$ws.Range("G14").Value = 0.6217654214631273

# Row 15 (HP_neutral / cornstover)
$ws.Range("D15").Value = 0.625825761106789
$ws.Range("E15").Value = 0.6088896977205107
$ws.Range("F15").Value = 0.6092256641852031
$ws.Range("G15").Value = 0.6094382465732844

# Row 16 (HP_neutral / sugarcane)
$ws.Range("D16").Value = 0.6258731801449422
$ws.Range("E16").Value = 0.6151550372018855
$ws.Range("F16").Value = 0.6157126652832703
$ws.Range("G16").Value = 0.6158215982194363

# Row 17 (HP_neutral / corn)
$ws.Range("D17").Value = 0.6263624387181024
$ws.Range("E17").Value = 0.6160755651237659
$ws.Range("F17").Value = 0.6161684793177287
$ws.Range("G17").Value = 0.6162814775809038
